$wb = $excel.ActiveWorkbook
$old = $wb.ActiveSheet
$oldName = $old.Name

# Duplicate the original sheet (keeps its full formatting/xr metadata) so the
# new sheet carries on the workbook's sheet-numbering, then rename it and
# drop the original - effectively "replacing" Sheet1 with ValidLogin.
$old.Copy($null, $old)
$copyName = $oldName + " (2)"
$new = $wb.Worksheets.Item($copyName)
$new.Name = "ValidLogin"
$wb.Worksheets.Item($oldName).Delete()

$ws = $wb.Worksheets.Item("ValidLogin")
$ws.Activate()

# Lay out the valid login credentials: headers in row 1, values in row 2.
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Match the saved view state: zoomed to 175% with B3 selected.
$excel.ActiveWindow.Zoom = 175
$ws.Range("B3").Select()
